# Update TPM-derived ligand-receptor edge statistics (Sema3c-Nrp2)
# per the new TPM-based recomputation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.034649
$ws.Range("H2").Value = 3.103947
$ws.Range("I2").Value = 0.02307585235784855
$ws.Range("J2").Value = 0.02307585235784855
$ws.Range("M2").Value = 70.23436
$ws.Range("N2").Value = 210.70308
$ws.Range("O2").Value = 0.7023186840741513
$ws.Range("P2").Value = 0.7023186840741513
$ws.Range("Q2").Value = 72.66791033963999
$ws.Range("R2").Value = 654.0111930567599
$ws.Range("S2").Value = 0.0162066022618536
$ws.Range("T2").Value = 0.0162066022618536

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.034649
$ws.Range("H3").Value = 3.103947
$ws.Range("I3").Value = 0.02307585235784855
$ws.Range("J3").Value = 0.02307585235784855
$ws.Range("O3").Value = 0.06551129587759326
$ws.Range("P3").Value = 0.06551129587759325
$ws.Range("Q3").Value = 6.778360141937999
$ws.Range("R3").Value = 61.005241277442
$ws.Range("S3").Value = 0.001511728991442675
$ws.Range("T3").Value = 0.001511728991442674

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.034649
$ws.Range("H4").Value = 3.103947
$ws.Range("I4").Value = 0.02307585235784855
$ws.Range("J4").Value = 0.02307585235784855
$ws.Range("M4").Value = 4.268944666666666
$ws.Range("N4").Value = 12.806834
$ws.Range("O4").Value = 0.04268793224112385
$ws.Range("P4").Value = 0.04268793224112385
$ws.Range("Q4").Value = 4.416859330421999
$ws.Range("R4").Value = 39.75173397379799
$ws.Range("S4").Value = 0.0009850604218580171
$ws.Range("T4").Value = 0.0009850604218580171

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.034649
$ws.Range("H5").Value = 3.103947
$ws.Range("I5").Value = 0.02307585235784855
$ws.Range("J5").Value = 0.02307585235784855
$ws.Range("M5").Value = 18.948881
$ws.Range("N5").Value = 56.846643
$ws.Range("O5").Value = 0.1894820878071316
$ws.Range("P5").Value = 0.1894820878071315
$ws.Range("Q5").Value = 19.605440777769
$ws.Range("R5").Value = 176.448966999921
$ws.Range("S5").Value = 0.004372460682694264
$ws.Range("T5").Value = 0.004372460682694263

# Row 6
$ws.Range("I6").Value = 0.9141262989281937
$ws.Range("J6").Value = 0.9141262989281937
$ws.Range("M6").Value = 70.23436
$ws.Range("N6").Value = 210.70308
$ws.Range("O6").Value = 0.7023186840741513
$ws.Range("P6").Value = 0.7023186840741513
$ws.Range("Q6").Value = 2878.664974081773
$ws.Range("R6").Value = 25907.98476673596
$ws.Range("S6").Value = 0.6420079793408233
$ws.Range("T6").Value = 0.6420079793408233

# Row 7
$ws.Range("I7").Value = 0.9141262989281937
$ws.Range("J7").Value = 0.9141262989281937
$ws.Range("O7").Value = 0.06551129587759326
$ws.Range("P7").Value = 0.06551129587759325
$ws.Range("S7").Value = 0.05988559843857416
$ws.Range("T7").Value = 0.05988559843857415

# Row 8
$ws.Range("I8").Value = 0.9141262989281937
$ws.Range("J8").Value = 0.9141262989281937
$ws.Range("M8").Value = 4.268944666666666
$ws.Range("N8").Value = 12.806834
$ws.Range("O8").Value = 0.04268793224112385
$ws.Range("P8").Value = 0.04268793224112385
$ws.Range("Q8").Value = 174.9693666778842
$ws.Range("R8").Value = 1574.724300100958
$ws.Range("S8").Value = 0.03902216150847607
$ws.Range("T8").Value = 0.03902216150847607

# Row 9
$ws.Range("I9").Value = 0.9141262989281937
$ws.Range("J9").Value = 0.9141262989281937
$ws.Range("M9").Value = 18.948881
$ws.Range("N9").Value = 56.846643
$ws.Range("O9").Value = 0.1894820878071316
$ws.Range("P9").Value = 0.1894820878071315
$ws.Range("Q9").Value = 776.6494922534156
$ws.Range("R9").Value = 6989.845430280741
$ws.Range("S9").Value = 0.1732105596403202
$ws.Range("T9").Value = 0.1732105596403202

# Row 10
$ws.Range("G10").Value = 2.775347666666667
$ws.Range("H10").Value = 8.326043
$ws.Range("I10").Value = 0.06189878209682655
$ws.Range("J10").Value = 0.06189878209682655
$ws.Range("M10").Value = 70.23436
$ws.Range("N10").Value = 210.70308
$ws.Range("O10").Value = 0.7023186840741513
$ws.Range("P10").Value = 0.7023186840741513
$ws.Range("Q10").Value = 194.9247671458266
$ws.Range("R10").Value = 1754.32290431244
$ws.Range("S10").Value = 0.04347267118803586
$ws.Range("T10").Value = 0.04347267118803586

# Row 11
$ws.Range("G11").Value = 2.775347666666667
$ws.Range("H11").Value = 8.326043
$ws.Range("I11").Value = 0.06189878209682655
$ws.Range("J11").Value = 0.06189878209682655
$ws.Range("O11").Value = 0.06551129587759326
$ws.Range("P11").Value = 0.06551129587759325
$ws.Range("Q11").Value = 18.18230724018867
$ws.Range("R11").Value = 163.640765161698
$ws.Range("S11").Value = 0.004055069428407877
$ws.Range("T11").Value = 0.004055069428407876

# Row 12
$ws.Range("G12").Value = 2.775347666666667
$ws.Range("H12").Value = 8.326043
$ws.Range("I12").Value = 0.06189878209682655
$ws.Range("J12").Value = 0.06189878209682655
$ws.Range("M12").Value = 4.268944666666666
$ws.Range("N12").Value = 12.806834
$ws.Range("O12").Value = 0.04268793224112385
$ws.Range("P12").Value = 0.04268793224112385
$ws.Range("Q12").Value = 11.84780561976244
$ws.Range("R12").Value = 106.630250577862
$ws.Range("S12").Value = 0.002642331015957422
$ws.Range("T12").Value = 0.002642331015957422

# Row 13
$ws.Range("G13").Value = 2.775347666666667
$ws.Range("H13").Value = 8.326043
$ws.Range("I13").Value = 0.06189878209682655
$ws.Range("J13").Value = 0.06189878209682655
$ws.Range("M13").Value = 18.948881
$ws.Range("N13").Value = 56.846643
$ws.Range("O13").Value = 0.1894820878071316
$ws.Range("P13").Value = 0.1894820878071315
$ws.Range("Q13").Value = 52.58973266929434
$ws.Range("R13").Value = 473.307594023649
$ws.Range("S13").Value = 0.01172871046442539
$ws.Range("T13").Value = 0.01172871046442539

# Row 14
$ws.Range("G14").Value = 0.04031133333333333
$ws.Range("H14").Value = 0.120934
$ws.Range("I14").Value = 0.0008990666171310454
$ws.Range("J14").Value = 0.0008990666171310454
$ws.Range("M14").Value = 70.23436
$ws.Range("N14").Value = 210.70308
$ws.Range("O14").Value = 0.7023186840741513
$ws.Range("P14").Value = 0.7023186840741513
$ws.Range("Q14").Value = 2.831240697413333
$ws.Range("R14").Value = 25.48116627672
$ws.Range("S14").Value = 0.0006314312834384746
$ws.Range("T14").Value = 0.0006314312834384746

# Row 15
$ws.Range("G15").Value = 0.04031133333333333
$ws.Range("H15").Value = 0.120934
$ws.Range("I15").Value = 0.0008990666171310454
$ws.Range("J15").Value = 0.0008990666171310454
$ws.Range("O15").Value = 0.06551129587759326
$ws.Range("P15").Value = 0.06551129587759325
$ws.Range("Q15").Value = 0.2640941373693333
$ws.Range("R15").Value = 2.376847236324
$ws.Range("S15").Value = 0.00005889901916853877
$ws.Range("T15").Value = 0.00005889901916853876

# Row 16
$ws.Range("G16").Value = 0.04031133333333333
$ws.Range("H16").Value = 0.120934
$ws.Range("I16").Value = 0.0008990666171310454
$ws.Range("J16").Value = 0.0008990666171310454
$ws.Range("M16").Value = 4.268944666666666
$ws.Range("N16").Value = 12.806834
$ws.Range("O16").Value = 0.04268793224112385
$ws.Range("P16").Value = 0.04268793224112385
$ws.Range("Q16").Value = 0.1720868514395555
$ws.Range("R16").Value = 1.548781662956
$ws.Range("S16").Value = 0.00003837929483234651
$ws.Range("T16").Value = 0.00003837929483234651

# Row 17
$ws.Range("G17").Value = 0.04031133333333333
$ws.Range("H17").Value = 0.120934
$ws.Range("I17").Value = 0.0008990666171310454
$ws.Range("J17").Value = 0.0008990666171310454
$ws.Range("M17").Value = 18.948881
$ws.Range("N17").Value = 56.846643
$ws.Range("O17").Value = 0.1894820878071316
$ws.Range("P17").Value = 0.1894820878071315
$ws.Range("Q17").Value = 0.7638546582846666
$ws.Range("R17").Value = 6.874691924562
$ws.Range("S17").Value = 0.0001703570196916855
$ws.Range("T17").Value = 0.0001703570196916855
